# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.255.82'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '3.705.80'
$ws.Range('E3').Value = '  +1.08%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '237.79'
$ws.Range('E5').Value = '  -0.74%  '
$ws.Range('E6').Value = '  +3.94%  '
$ws.Range('D7').Value = '660.54'
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('E8').Value = '  +1.47%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  -1.33%  '
$ws.Range('D11').Value = '3.704.51'
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').Value = '0.0000322'
$ws.Range('E12').Value = '  +20.95%  '
$ws.Range('D13').Value = '44.51'
$ws.Range('E13').Value = '  -2.80%  '
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('D15').Value = '6.89'
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').Value = '4.398.92'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('D17').Value = '97.110.98'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('E18').Value = '  +16.88%  '
$ws.Range('D19').Value = '3.690.66'
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('D20').Value = '13.10'
$ws.Range('E20').Value = '  +2.37%  '
$ws.Range('D21').Value = '18.79'
$ws.Range('E21').Value = '  +1.08%  '
$ws.Range('E22').Value = '  -3.20%  '
$ws.Range('D23').Value = '523.88'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = '3.45'
$ws.Range('E24').Value = '  +0.69%  '
$ws.Range('E25').Value = '  +7.50%  '
$ws.Range('D26').Value = '6.92'
$ws.Range('E26').Value = '  -2.35%  '
$ws.Range('D27').Value = '103.07'
$ws.Range('E27').Value = '  +1.65%  '
$ws.Range('E28').Value = '  +15.65%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '3.908.81'
$ws.Range('E29').Value = '  +1.17%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = '13.52'
$ws.Range('E30').Value = '  +2.96%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '12.97'
$ws.Range('E31').Value = '  +4.08%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '3.05'
$ws.Range('E32').Value = '  +0.98%  '
$ws.Range('B33').Value = 'Dai'
$ws.Range('C33').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D33').Value = '0.999'
$ws.Range('E33').Value = '  +0.15%  '
$ws.Range('B34').Value = 'Cronos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D34').Value = '0.193'
$ws.Range('E34').Value = '  +4.39%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = '1.86'
$ws.Range('E35').Value = '  -1.86%  '
$ws.Range('B36').Value = 'Binance-PegBSC-USD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D36').Value = '0.993'
$ws.Range('E36').Value = '  -1.07%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '653.99'
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').Value = '0.597'
$ws.Range('E39').Value = '  +1.54%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D40').Value = '8.85'
$ws.Range('E40').Value = '  +0.60%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.168'
$ws.Range('E42').Value = '  +5.43%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '6.89'
$ws.Range('E43').Value = '  +6.29%  '
$ws.Range('B44').Value = 'ImmutableX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D44').Value = '2.06'
$ws.Range('E44').Value = '  +4.86%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = '0.491'
$ws.Range('E45').Value = '  +11.42%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '40.30'
$ws.Range('E46').Value = '  +1.44%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '0.973'
$ws.Range('E47').Value = '  +1.89%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0463'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '2.45'
$ws.Range('E49').Value = '  +5.82%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '8.81'
$ws.Range('E50').Value = '  +3.05%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '23.63'
$ws.Range('E51').Value = '  -0.13%  '
